# GRNmap dCIN5 workbook edit:
#   - "production_rates" and "degradation_rates" sheets: the B2:B16 values were
#     stored with the wrong sign (all negative) - flip them to positive.
#   - "dcin5_log2_expression" sheet: the time-point headings (row 1, F:M) were
#     all "15" - correct them to 30 (F:I) and 60 (J:M) so the sheet reflects the
#     15/30/60-minute time course instead of all-15s.
#   - Update the remembered cell selections on several sheets to match where the
#     author last left the cursor, and move the "active" tab from
#     optimization_parameters to dhap4_log2_expression.

$wb = $excel.ActiveWorkbook

# --- production_rates: negate production rate column (B2:B16) ---
$wsProd = $wb.Worksheets.Item("production_rates")
for ($r = 2; $r -le 16; $r++) {
    $cell = $wsProd.Cells.Item($r, 2)
    $cell.Value = -($cell.Value2)
}
$wsProd.Range("B18").Select()

# --- degradation_rates: negate degradation rate column (B2:B16) ---
$wsDeg = $wb.Worksheets.Item("degradation_rates")
for ($r = 2; $r -le 16; $r++) {
    $cell = $wsDeg.Cells.Item($r, 2)
    $cell.Value = -($cell.Value2)
}
$wsDeg.Range("A2:A16").Select()

# --- wt_log2_expression: selection only ---
$wsWt = $wb.Worksheets.Item("wt_log2_expression")
$wsWt.Range("A2:A16").Select()

# --- dcin5_log2_expression: fix the time headings in row 1 ---
$wsCin5 = $wb.Worksheets.Item("dcin5_log2_expression")
$wsCin5.Range("F1:I1").Value = 30
$wsCin5.Range("J1:M1").Value = 60
$wsCin5.Range("L14").Select()

# --- dgln3_log2_expression: selection only ---
$wsGln3 = $wb.Worksheets.Item("dgln3_log2_expression")
$wsGln3.Range("A2:A16").Select()

# --- dhap4_log2_expression: selection + becomes the active/visible tab ---
# Activate this one last so it ends up as the workbook's active sheet
# (matches activeTab moving from optimization_parameters to this sheet).
$wsHap4 = $wb.Worksheets.Item("dhap4_log2_expression")
$wsHap4.Activate()
$wsHap4.Range("D18").Select()
